$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.594.65'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.203.22'
$ws.Range("E3").Value = '  -1.94%  '
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = '  -0.08%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.29'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -1.16%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.618'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -3.07%  '
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.94'
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = '  -6.14%  '
$ws.Range("E8").Value = '  -0.03%  '
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.400'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  -2.21%  '
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.43'
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = '  -3.38%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("E12").Value = '  -1.21%  '
$ws.Range("D13").Value = '2.528.83'
$ws.Range("E13").Value = '  -2.15%  '
$ws.Range("E14").Value = '  -4.86%  '
$ws.Range("E15").Value = '  -1.43%  '
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.63'
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = '  -0.48%  '
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.792'
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = '  -4.06%  '
$ws.Range("D18").Value = '2.235.17'
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").Value = '41.502.02'
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").Value = '0.0₃0898'
$ws.Range("E20").Value = '  -2.66%  '
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.85'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  -2.33%  '
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  -2.12%  '
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.48'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  -3.49%  '
$ws.Range("E24").Value = '  -0.15%  '
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  +1.83%  '
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.35'
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  -2.09%  '
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.63'
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = '  -2.69%  '
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.79'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  -2.54%  '
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.138'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  -4.86%  '
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.68'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  -3.76%  '
$ws.Range("E31").Value = '  -2.37%  '
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.52'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  -10.23%  '
$ws.Range("E33").Value = '  -3.09%  '
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.98'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  -1.59%  '
$ws.Range("E35").Value = '  -2.49%  '
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0646'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  +1.86%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.43'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  -7.16%  '
$ws.Range("E38").Value = '  -3.70%  '
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.56'
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  -7.08%  '
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  -0.24%  '
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.000235'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  -13.37%  '
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0236'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  -1.87%  '
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.51'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  -3.91%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.20'
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = '  -2.44%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0954'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  +1.02%  '
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '97.13'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  -5.48%  '
$ws.Range("D47").Value = '1.464.24'
$ws.Range("E47").Value = '  -3.11%  '
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.33'
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  -11.83%  '
$ws.Range("E49").Value = '  -7.37%  '
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.75'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  -1.29%  '
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.07'
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  -4.67%  '
